# Auto-generated: update cached market/profit figures across all 8 class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match refreshed Market Board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1544.7084
$ws.Range("J17").Value = 1696.2106
$ws.Range("L17").Value = 5088.6318
$ws.Range("N17").Value = -5424.6318
$ws.Range("H33").Value = 226.27777
$ws.Range("I33").Value = 200.26666
$ws.Range("J33").Value = 356.33334
$ws.Range("K33").Value = 200.26666
$ws.Range("L33").Value = 356.33334
$ws.Range("M33").Value = 28.73334
$ws.Range("N33").Value = -814.33334
$ws.Range("H40").Value = 5137.2104
$ws.Range("J40").Value = 3124
$ws.Range("L40").Value = 3124
$ws.Range("N40").Value = -3474
$ws.Range("H113").Value = 3879.8
$ws.Range("I113").Value = 3599.75
$ws.Range("K113").Value = 3599.75
$ws.Range("M113").Value = -345.75
$ws.Range("H125").Value = 1032.381
$ws.Range("J125").Value = 849.4
$ws.Range("L125").Value = 7644.599999999999
$ws.Range("N125").Value = -12564.6
$ws.Range("H130").Value = 125560.5
$ws.Range("J130").Value = 125560.5
$ws.Range("L130").Value = 125560.5
$ws.Range("N130").Value = -135600.5
$ws.Range("H132").Value = 53210.74
$ws.Range("I132").Value = 61717.9
$ws.Range("J132").Value = 5813.7144
$ws.Range("K132").Value = 185153.7
$ws.Range("L132").Value = 17441.1432
$ws.Range("M132").Value = -182623.7
$ws.Range("N132").Value = -22501.1432
$ws.Range("H135").Value = 971.06976
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 24392330
$ws.Range("I137").Value = 34484416
$ws.Range("K137").Value = 103453248
$ws.Range("M137").Value = -103450698

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8477834
$ws.Range("I32").Value = 9806202
$ws.Range("K32").Value = 9806202
$ws.Range("M32").Value = -9805915
$ws.Range("H74").Value = 4466806
$ws.Range("I74").Value = 5683302.5
$ws.Range("J74").Value = 6319.6665
$ws.Range("K74").Value = 5683302.5
$ws.Range("L74").Value = 6319.6665
$ws.Range("M74").Value = -5682428.5
$ws.Range("N74").Value = -8067.6665
$ws.Range("H77").Value = 4466806
$ws.Range("I77").Value = 5683302.5
$ws.Range("J77").Value = 6319.6665
$ws.Range("K77").Value = 28416512.5
$ws.Range("L77").Value = 31598.3325
$ws.Range("M77").Value = -28412144.5
$ws.Range("N77").Value = -40334.3325
$ws.Range("H132").Value = 519417.44
$ws.Range("I132").Value = 531119.75
$ws.Range("K132").Value = 1593359.25
$ws.Range("M132").Value = -1590829.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2712.3845
$ws.Range("I107").Value = 2686.7
$ws.Range("K107").Value = 2686.7
$ws.Range("M107").Value = -766.6999999999998
$ws.Range("H134").Value = 406029.47
$ws.Range("I134").Value = 555404.75
$ws.Range("J134").Value = 4583.4375
$ws.Range("K134").Value = 1666214.25
$ws.Range("L134").Value = 13750.3125
$ws.Range("M134").Value = -1663679.25
$ws.Range("N134").Value = -18820.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 123543.63
$ws.Range("I31").Value = 155591.47
$ws.Range("J31").Value = 48765.332
$ws.Range("K31").Value = 155591.47
$ws.Range("L31").Value = 48765.332
$ws.Range("M31").Value = -155296.47
$ws.Range("N31").Value = -49355.332
$ws.Range("H34").Value = 123543.63
$ws.Range("I34").Value = 155591.47
$ws.Range("J34").Value = 48765.332
$ws.Range("K34").Value = 155591.47
$ws.Range("L34").Value = 48765.332
$ws.Range("M34").Value = -155389.47
$ws.Range("N34").Value = -49169.332
$ws.Range("H58").Value = 327470.78
$ws.Range("I58").Value = 476538.22
$ws.Range("K58").Value = 476538.22
$ws.Range("M58").Value = -476335.22
$ws.Range("H59").Value = 63888.332
$ws.Range("J59").Value = 63888.332
$ws.Range("L59").Value = 63888.332
$ws.Range("N59").Value = -66178.33199999999
$ws.Range("H60").Value = 37083.168
$ws.Range("I60").Value = 20625
$ws.Range("J60").Value = 69999.5
$ws.Range("K60").Value = 20625
$ws.Range("L60").Value = 69999.5
$ws.Range("M60").Value = -20114
$ws.Range("N60").Value = -71021.5
$ws.Range("H81").Value = 102773.336
$ws.Range("J81").Value = 102773.336
$ws.Range("L81").Value = 102773.336
$ws.Range("N81").Value = -104769.336
$ws.Range("H84").Value = 102773.336
$ws.Range("J84").Value = 102773.336
$ws.Range("L84").Value = 308320.008
$ws.Range("N84").Value = -318304.008
$ws.Range("H132").Value = 44656644
$ws.Range("I132").Value = 38476344
$ws.Range("J132").Value = 125000500
$ws.Range("K132").Value = 115429032
$ws.Range("L132").Value = 375001500
$ws.Range("M132").Value = -115426502
$ws.Range("N132").Value = -375006560
$ws.Range("H136").Value = 327470.78
$ws.Range("I136").Value = 476538.22
$ws.Range("K136").Value = 1429614.66
$ws.Range("M136").Value = -1427064.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 2064.5557
$ws.Range("I24").Value = 299.66666
$ws.Range("J24").Value = 2947
$ws.Range("K24").Value = 898.9999799999999
$ws.Range("L24").Value = 8841
$ws.Range("M24").Value = -668.9999799999999
$ws.Range("N24").Value = -9301
$ws.Range("H33").Value = 147.54546
$ws.Range("I33").Value = 141.8421
$ws.Range("J33").Value = 183.66667
$ws.Range("K33").Value = 851.0526
$ws.Range("L33").Value = 1102.00002
$ws.Range("M33").Value = -568.0526
$ws.Range("N33").Value = -1668.00002
$ws.Range("H38").Value = 111.71429
$ws.Range("I38").Value = 33.6
$ws.Range("J38").Value = 155.11111
$ws.Range("K38").Value = 100.8
$ws.Range("L38").Value = 465.33333
$ws.Range("M38").Value = 246.2
$ws.Range("N38").Value = -1159.33333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1799.8572
$ws.Range("I3").Value = 319
$ws.Range("J3").Value = 5502
$ws.Range("K3").Value = 319
$ws.Range("L3").Value = 5502
$ws.Range("M3").Value = -203
$ws.Range("N3").Value = -5734
$ws.Range("H132").Value = 261733.94
$ws.Range("I132").Value = 276010.94
$ws.Range("J132").Value = 52338
$ws.Range("K132").Value = 828032.8200000001
$ws.Range("L132").Value = 157014
$ws.Range("M132").Value = -825502.8200000001
$ws.Range("N132").Value = -162074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1991.6666
$ws.Range("I22").Value = 1711.1111
$ws.Range("J22").Value = 2833.3333
$ws.Range("K22").Value = 1711.1111
$ws.Range("L22").Value = 2833.3333
$ws.Range("M22").Value = -1416.1111
$ws.Range("N22").Value = -3423.3333
$ws.Range("H27").Value = 1991.6666
$ws.Range("I27").Value = 1711.1111
$ws.Range("J27").Value = 2833.3333
$ws.Range("K27").Value = 1711.1111
$ws.Range("L27").Value = 2833.3333
$ws.Range("M27").Value = -1604.1111
$ws.Range("N27").Value = -3047.3333
$ws.Range("H46").Value = 5053.294
$ws.Range("I46").Value = 5911.2144
$ws.Range("J46").Value = 1049.6666
$ws.Range("K46").Value = 5911.2144
$ws.Range("L46").Value = 1049.6666
$ws.Range("M46").Value = -5723.2144
$ws.Range("N46").Value = -1425.6666
$ws.Range("H131").Value = 20326
$ws.Range("J131").Value = 20326
$ws.Range("L131").Value = 20326
$ws.Range("N131").Value = -30406
$ws.Range("H132").Value = 612554.1
$ws.Range("I132").Value = 830448.1
$ws.Range("K132").Value = 2491344.3
$ws.Range("M132").Value = -2488814.3
$ws.Range("H136").Value = 39985.363
$ws.Range("I136").Value = 2216.7932
$ws.Range("K136").Value = 6650.3796
$ws.Range("M136").Value = -4100.3796

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3250
$ws.Range("J14").Value = 3250
$ws.Range("L14").Value = 3250
$ws.Range("N14").Value = -3586
$ws.Range("H132").Value = 9334062
$ws.Range("I132").Value = 9788919
$ws.Range("K132").Value = 29366757
$ws.Range("M132").Value = -29364227
$ws.Range("H136").Value = 10814929
$ws.Range("I136").Value = 13719213
$ws.Range("K136").Value = 41157639
$ws.Range("M136").Value = -41155089

